$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Cells.Item(12, 8).Value = 111.78571
$ws.Cells.Item(12, 9).Value = 97.30768999999999
$ws.Cells.Item(12, 10).Value = 300
$ws.Cells.Item(12, 11).Value = 97.30768999999999
$ws.Cells.Item(12, 12).Value = 300
$ws.Cells.Item(12, 13).Value = 72.69231000000001
$ws.Cells.Item(12, 14).Value = -640
$ws.Cells.Item(17, 8).Value = 653.3333
$ws.Cells.Item(17, 10).Value = 685.7143
$ws.Cells.Item(17, 12).Value = 2057.1429
$ws.Cells.Item(17, 14).Value = -2393.1429
$ws.Cells.Item(107, 8).Value = 1041.0435
$ws.Cells.Item(107, 9).Value = 1087.8096
$ws.Cells.Item(107, 10).Value = 550
$ws.Cells.Item(107, 11).Value = 1087.8096
$ws.Cells.Item(107, 12).Value = 550
$ws.Cells.Item(107, 13).Value = 832.1904
$ws.Cells.Item(107, 14).Value = -4390
$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Cells.Item(4, 8).Value = 1000
$ws.Cells.Item(4, 10).Value = 1000
$ws.Cells.Item(4, 12).Value = 1000
$ws.Cells.Item(4, 14).Value = -1232
$ws.Cells.Item(5, 8).Value = 100000000
$ws.Cells.Item(5, 9).Value = 100000000
$ws.Cells.Item(5, 11).Value = 100000000
$ws.Cells.Item(5, 13).Value = -99999888
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 13).Value = $null
$ws.Cells.Item(9, 8).Value = 36605.4
$ws.Cells.Item(9, 9).Value = 1500
$ws.Cells.Item(9, 10).Value = 60009
$ws.Cells.Item(9, 11).Value = 1500
$ws.Cells.Item(9, 12).Value = 60009
$ws.Cells.Item(9, 13).Value = -1330
$ws.Cells.Item(9, 14).Value = -60349
$ws.Cells.Item(20, 8).Value = 36605.4
$ws.Cells.Item(20, 9).Value = 1500
$ws.Cells.Item(20, 10).Value = 60009
$ws.Cells.Item(20, 11).Value = 1500
$ws.Cells.Item(20, 12).Value = 60009
$ws.Cells.Item(20, 13).Value = -1230
$ws.Cells.Item(20, 14).Value = -60549
$ws.Cells.Item(23, 8).Value = 38341.4
$ws.Cells.Item(23, 10).Value = 41926.75
$ws.Cells.Item(23, 12).Value = 41926.75
$ws.Cells.Item(23, 14).Value = -42444.75
$ws.Cells.Item(37, 8).Value = 12620.083
$ws.Cells.Item(37, 10).Value = 16648.715
$ws.Cells.Item(37, 12).Value = 16648.715
$ws.Cells.Item(37, 14).Value = -17194.715
$ws.Cells.Item(44, 8).Value = 21714.285
$ws.Cells.Item(44, 10).Value = 21714.285
$ws.Cells.Item(44, 12).Value = 21714.285
$ws.Cells.Item(44, 14).Value = -22690.285
$ws.Cells.Item(55, 8).Value = 51167.168
$ws.Cells.Item(55, 10).Value = 51167.168
$ws.Cells.Item(55, 12).Value = 51167.168
$ws.Cells.Item(55, 14).Value = -51797.168
$ws.Cells.Item(63, 8).Value = 1650.5
$ws.Cells.Item(63, 9).Value = 1412.8125
$ws.Cells.Item(63, 10).Value = 2284.3333
$ws.Cells.Item(63, 11).Value = 1412.8125
$ws.Cells.Item(63, 12).Value = 2284.3333
$ws.Cells.Item(63, 13).Value = -726.8125
$ws.Cells.Item(63, 14).Value = -3656.3333
$ws.Cells.Item(66, 8).Value = 1650.5
$ws.Cells.Item(66, 9).Value = 1412.8125
$ws.Cells.Item(66, 10).Value = 2284.3333
$ws.Cells.Item(66, 11).Value = 7064.0625
$ws.Cells.Item(66, 12).Value = 11421.6665
$ws.Cells.Item(66, 13).Value = -3632.0625
$ws.Cells.Item(66, 14).Value = -18285.6665
$ws.Cells.Item(74, 8).Value = 36509256
$ws.Cells.Item(74, 9).Value = 28572116
$ws.Cells.Item(74, 11).Value = 28572116
$ws.Cells.Item(74, 13).Value = -28571242
$ws.Cells.Item(77, 8).Value = 36509256
$ws.Cells.Item(77, 9).Value = 28572116
$ws.Cells.Item(77, 11).Value = 142860580
$ws.Cells.Item(77, 13).Value = -142856212
$ws.Cells.Item(80, 8).Value = 23375.25
$ws.Cells.Item(80, 10).Value = 23375.25
$ws.Cells.Item(80, 12).Value = 23375.25
$ws.Cells.Item(80, 14).Value = -25371.25
$ws.Cells.Item(83, 8).Value = 23375.25
$ws.Cells.Item(83, 10).Value = 23375.25
$ws.Cells.Item(83, 12).Value = 70125.75
$ws.Cells.Item(83, 14).Value = -80109.75
$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Cells.Item(4, 8).Value = 100000000
$ws.Cells.Item(4, 9).Value = 100000000
$ws.Cells.Item(4, 11).Value = 100000000
$ws.Cells.Item(4, 13).Value = -99999885
$ws.Cells.Item(15, 8).Value = 27669
$ws.Cells.Item(15, 9).Value = 11500
$ws.Cells.Item(15, 10).Value = 60007
$ws.Cells.Item(15, 11).Value = 11500
$ws.Cells.Item(15, 12).Value = 60007
$ws.Cells.Item(15, 13).Value = -11273
$ws.Cells.Item(15, 14).Value = -60461
$ws.Cells.Item(19, 8).Value = 10000
$ws.Cells.Item(19, 9).Value = 10000
$ws.Cells.Item(19, 11).Value = 10000
$ws.Cells.Item(19, 13).Value = -9827
$ws.Cells.Item(22, 8).Value = 393.32144
$ws.Cells.Item(22, 9).Value = 446.04544
$ws.Cells.Item(22, 10).Value = 200
$ws.Cells.Item(22, 11).Value = 446.04544
$ws.Cells.Item(22, 12).Value = 200
$ws.Cells.Item(22, 13).Value = -273.04544
$ws.Cells.Item(22, 14).Value = -546
$ws.Cells.Item(35, 8).Value = 23990
$ws.Cells.Item(35, 10).Value = 23990
$ws.Cells.Item(35, 12).Value = 23990
$ws.Cells.Item(35, 14).Value = -24610
$ws.Cells.Item(82, 8).Value = 17147.5
$ws.Cells.Item(82, 9).Value = 1833.8572
$ws.Cells.Item(82, 10).Value = 25393.309
$ws.Cells.Item(82, 11).Value = 1833.8572
$ws.Cells.Item(82, 12).Value = 25393.309
$ws.Cells.Item(82, 13).Value = -1450.8572
$ws.Cells.Item(82, 14).Value = -26159.309
$ws.Cells.Item(85, 8).Value = 17147.5
$ws.Cells.Item(85, 9).Value = 1833.8572
$ws.Cells.Item(85, 10).Value = 25393.309
$ws.Cells.Item(85, 11).Value = 1833.8572
$ws.Cells.Item(85, 12).Value = 25393.309
$ws.Cells.Item(85, 13).Value = -507.8571999999999
$ws.Cells.Item(85, 14).Value = -28045.309
$ws.Cells.Item(132, 8).Value = 39380
$ws.Cells.Item(132, 10).Value = 39380
$ws.Cells.Item(132, 12).Value = 39380
$ws.Cells.Item(132, 14).Value = -49500
$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Cells.Item(7, 8).Value = 20408258
$ws.Cells.Item(7, 9).Value = 47619096
$ws.Cells.Item(7, 10).Value = 130
$ws.Cells.Item(7, 11).Value = 47619096
$ws.Cells.Item(7, 12).Value = 130
$ws.Cells.Item(7, 13).Value = -47618983
$ws.Cells.Item(7, 14).Value = -356
$ws.Cells.Item(31, 8).Value = 1160057
$ws.Cells.Item(31, 9).Value = 954.97675
$ws.Cells.Item(31, 11).Value = 954.97675
$ws.Cells.Item(31, 13).Value = -659.97675
$ws.Cells.Item(34, 8).Value = 1160057
$ws.Cells.Item(34, 9).Value = 954.97675
$ws.Cells.Item(34, 11).Value = 954.97675
$ws.Cells.Item(34, 13).Value = -752.97675
$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Cells.Item(46, 8).Value = 2152.5386
$ws.Cells.Item(46, 10).Value = 2312.0908
$ws.Cells.Item(46, 12).Value = 6936.2724
$ws.Cells.Item(46, 14).Value = -7118.2724
$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Cells.Item(46, 8).Value = 1112.1333
$ws.Cells.Item(46, 9).Value = 1020.2857
$ws.Cells.Item(46, 10).Value = 1192.5
$ws.Cells.Item(46, 11).Value = 1020.2857
$ws.Cells.Item(46, 12).Value = 1192.5
$ws.Cells.Item(46, 13).Value = -832.2857
$ws.Cells.Item(46, 14).Value = -1568.5
$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Cells.Item(132, 8).Value = 822664.9399999999
$ws.Cells.Item(132, 9).Value = 1494.1875
$ws.Cells.Item(132, 10).Value = 4764284.5
$ws.Cells.Item(132, 11).Value = 4482.5625
$ws.Cells.Item(132, 12).Value = 14292853.5
$ws.Cells.Item(132, 13).Value = -1952.5625
$ws.Cells.Item(132, 14).Value = -14297913.5
